$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# Fix the C18 cell value: remove stray leading space/apostrophe artifact,
# store it as quote-prefixed text (so the shared string stays clean and
# the cell gets the quotePrefix style), matching the Latitude/Longitude
# text-number convention used throughout the sheet.
$ws.Range("C18").Value = "'-9.227705326284289"

# Update the active view: scroll so row 2 is the top-left visible row,
# and move the selection to D15.
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("D15").Select()
